$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 64
$arr64 = New-Object 'object[,]' 1,28
$arr64[0,0] = 7082624
$arr64[0,1] = "Chile Primera Division"
$arr64[0,2] = "Chile Primera Division"
$arr64[0,3] = 45183.79166666666
$arr64[0,4] = "Colo Colo"
$arr64[0,5] = "Deportes Copiapo"
$arr64[0,6] = 1
$arr64[0,7] = 1
$arr64[0,8] = "D"
$arr64[0,9] = 1.333
$arr64[0,10] = 5
$arr64[0,11] = 8
$arr64[0,12] = 1.45
$arr64[0,13] = 4.75
$arr64[0,14] = 7
$arr64[0,15] = -1.25
$arr64[0,16] = 1.925
$arr64[0,17] = 1.875
$arr64[0,18] = 3
$arr64[0,19] = 1.875
$arr64[0,20] = 1.925
$arr64[0,21] = -1
$arr64[0,22] = 3.75
$arr64[0,23] = -1
$arr64[0,24] = -1
$arr64[0,25] = 0.875
$arr64[0,26] = -1
$arr64[0,27] = 0.925
$ws.Range("B64:AC64").Value = $arr64

# Row 65
$arr65 = New-Object 'object[,]' 1,28
$arr65[0,0] = 7157967
$arr65[0,1] = "Chile Primera Division"
$arr65[0,2] = "Chile Primera Division"
$arr65[0,3] = 45183.79166666666
$arr65[0,4] = "Huachipato"
$arr65[0,5] = "Palestino"
$arr65[0,6] = 2
$arr65[0,7] = 2
$arr65[0,8] = "D"
$arr65[0,9] = 2.375
$arr65[0,10] = 3.2
$arr65[0,11] = 3
$arr65[0,12] = 2.75
$arr65[0,13] = 3.2
$arr65[0,14] = 2.7
$arr65[0,15] = 0
$arr65[0,16] = 1.925
$arr65[0,17] = 1.875
$arr65[0,18] = 2.5
$arr65[0,19] = 2
$arr65[0,20] = 1.8
$arr65[0,21] = -1
$arr65[0,22] = 2.2
$arr65[0,23] = -1
$arr65[0,24] = 0
$arr65[0,25] = -0
$arr65[0,26] = 1
$arr65[0,27] = -1
$ws.Range("B65:AC65").Value = $arr65

# Row 97
$arr97 = New-Object 'object[,]' 1,28
$arr97[0,0] = 7323253
$arr97[0,1] = "Chile Primera Division"
$arr97[0,2] = "Chile Primera Division"
$arr97[0,3] = 45242.83333333334
$arr97[0,4] = "Union Espanola"
$arr97[0,5] = "OHiggins"
$arr97[0,6] = 3
$arr97[0,7] = 3
$arr97[0,8] = "D"
$arr97[0,9] = 2
$arr97[0,10] = 3.4
$arr97[0,11] = 3.5
$arr97[0,12] = 2.1
$arr97[0,13] = 3.5
$arr97[0,14] = 3.75
$arr97[0,15] = -0.5
$arr97[0,16] = 2.025
$arr97[0,17] = 1.775
$arr97[0,18] = 2.5
$arr97[0,19] = 1.95
$arr97[0,20] = 1.85
$arr97[0,21] = -1
$arr97[0,22] = 2.5
$arr97[0,23] = -1
$arr97[0,24] = -1
$arr97[0,25] = 0.7749999999999999
$arr97[0,26] = 0.95
$arr97[0,27] = -1
$ws.Range("B97:AC97").Value = $arr97

# Row 98
$arr98 = New-Object 'object[,]' 1,28
$arr98[0,0] = 7323186
$arr98[0,1] = "Chile Primera Division"
$arr98[0,2] = "Chile Primera Division"
$arr98[0,3] = 45242.83333333334
$arr98[0,4] = "Coquimbo Unido"
$arr98[0,5] = "Deportes Copiapo"
$arr98[0,6] = 1
$arr98[0,7] = 0
$arr98[0,8] = "H"
$arr98[0,9] = 2
$arr98[0,10] = 3.4
$arr98[0,11] = 3.5
$arr98[0,12] = 1.727
$arr98[0,13] = 3.8
$arr98[0,14] = 4.75
$arr98[0,15] = -0.75
$arr98[0,16] = 1.9
$arr98[0,17] = 1.9
$arr98[0,18] = 2.75
$arr98[0,19] = 1.85
$arr98[0,20] = 1.95
$arr98[0,21] = 0.7270000000000001
$arr98[0,22] = -1
$arr98[0,23] = -1
$arr98[0,24] = 0.45
$arr98[0,25] = -0.5
$arr98[0,26] = -1
$arr98[0,27] = 0.95
$ws.Range("B98:AC98").Value = $arr98

# Row 102
$arr102 = New-Object 'object[,]' 1,28
$arr102[0,0] = 7494646
$arr102[0,1] = "Chile Primera Division"
$arr102[0,2] = "Chile Primera Division"
$arr102[0,3] = 45255.75
$arr102[0,4] = "OHiggins"
$arr102[0,5] = "Cobresal"
$arr102[0,6] = 0
$arr102[0,7] = 0
$arr102[0,8] = "D"
$arr102[0,9] = 3
$arr102[0,10] = 3.4
$arr102[0,11] = 2.3
$arr102[0,12] = 2.1
$arr102[0,13] = 3.5
$arr102[0,14] = 3.5
$arr102[0,15] = -0.25
$arr102[0,16] = 1.8
$arr102[0,17] = 2.05
$arr102[0,18] = 2.75
$arr102[0,19] = 1.975
$arr102[0,20] = 1.875
$arr102[0,21] = -1
$arr102[0,22] = 2.5
$arr102[0,23] = -1
$arr102[0,24] = -0.5
$arr102[0,25] = 0.5249999999999999
$arr102[0,26] = -1
$arr102[0,27] = 0.875
$ws.Range("B102:AC102").Value = $arr102

# Row 103
$arr103 = New-Object 'object[,]' 1,28
$arr103[0,0] = 7494647
$arr103[0,1] = "Chile Primera Division"
$arr103[0,2] = "Chile Primera Division"
$arr103[0,3] = 45255.75
$arr103[0,4] = "Huachipato"
$arr103[0,5] = "Universidad Catolica"
$arr103[0,6] = 1
$arr103[0,7] = 1
$arr103[0,8] = "D"
$arr103[0,9] = 2.2
$arr103[0,10] = 3.4
$arr103[0,11] = 3.2
$arr103[0,12] = 1.8
$arr103[0,13] = 3.6
$arr103[0,14] = 4.333
$arr103[0,15] = -0.75
$arr103[0,16] = 1.975
$arr103[0,17] = 1.875
$arr103[0,18] = 2.75
$arr103[0,19] = 1.975
$arr103[0,20] = 1.875
$arr103[0,21] = -1
$arr103[0,22] = 2.6
$arr103[0,23] = -1
$arr103[0,24] = -1
$arr103[0,25] = 0.875
$arr103[0,26] = -1
$arr103[0,27] = 0.875
$ws.Range("B103:AC103").Value = $arr103

# Row 137
$arr137 = New-Object 'object[,]' 1,28
$arr137[0,0] = 7723533
$arr137[0,1] = "Chile Primera Division"
$arr137[0,2] = "Chile Primera Division"
$arr137[0,3] = 45353.75
$arr137[0,4] = "OHiggins"
$arr137[0,5] = "Everton de Vina"
$arr137[0,6] = 2
$arr137[0,7] = 1
$arr137[0,8] = "H"
$arr137[0,9] = 3
$arr137[0,10] = 3.2
$arr137[0,11] = 2.375
$arr137[0,12] = 2.3
$arr137[0,13] = 3.1
$arr137[0,14] = 3.3
$arr137[0,15] = -0.25
$arr137[0,16] = 1.9
$arr137[0,17] = 1.9
$arr137[0,18] = 2.25
$arr137[0,19] = 1.95
$arr137[0,20] = 1.85
$arr137[0,21] = 1.3
$arr137[0,22] = -1
$arr137[0,23] = -1
$arr137[0,24] = 0.8999999999999999
$arr137[0,25] = -1
$arr137[0,26] = 0.95
$arr137[0,27] = -1
$ws.Range("B137:AC137").Value = $arr137

# Row 138
$arr138 = New-Object 'object[,]' 1,28
$arr138[0,0] = 7723528
$arr138[0,1] = "Chile Primera Division"
$arr138[0,2] = "Chile Primera Division"
$arr138[0,3] = 45353.75
$arr138[0,4] = "Palestino"
$arr138[0,5] = "Universidad Catolica"
$arr138[0,6] = 0
$arr138[0,7] = 2
$arr138[0,8] = "A"
$arr138[0,9] = 1.95
$arr138[0,10] = 3.6
$arr138[0,11] = 3.4
$arr138[0,12] = 2.375
$arr138[0,13] = 3.5
$arr138[0,14] = 2.9
$arr138[0,15] = 0
$arr138[0,16] = 1.8
$arr138[0,17] = 2.05
$arr138[0,18] = 2.25
$arr138[0,19] = 1.8
$arr138[0,20] = 2.05
$arr138[0,21] = -1
$arr138[0,22] = -1
$arr138[0,23] = 1.9
$arr138[0,24] = -1
$arr138[0,25] = 1.05
$arr138[0,26] = -0.5
$arr138[0,27] = 0.5249999999999999
$ws.Range("B138:AC138").Value = $arr138

# Row 177
$ws.Range("N177").Value = 2.625
$ws.Range("O177").Value = 3.2
$ws.Range("R177").Value = 1.875
$ws.Range("S177").Value = 1.975

# Row 180
$ws.Range("N180").Value = 2.75
$ws.Range("P180").Value = 2.7
$ws.Range("R180").Value = 1.925
$ws.Range("S180").Value = 1.925

# Row 181
$ws.Range("U181").Value = 1.875
$ws.Range("V181").Value = 1.975
